$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

$ws.Range("D2").Value = '37.433.59'
$ws.Range("E2").Value = '  +1.11%  '
$ws.Range("D3").Value = '2.033.08'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("E4").Value = '  +0.35%  '
Set-TextValue "D5" '229.57'
$ws.Range("E5").Value = '  +1.00%  '
$ws.Range("E6").Value = '  +0.80%  '
$ws.Range("E7").Value = '  +0.06%  '
Set-TextValue "D8" '56.43'
$ws.Range("E8").Value = '  +2.33%  '
Set-TextValue "D9" '0.381'
$ws.Range("E9").Value = '  -0.20%  '
Set-TextValue "D10" '0.0786'
$ws.Range("E10").Value = '  -1.04%  '
$ws.Range("D12").Value = '2.336.65'
$ws.Range("E12").Value = '  +0.56%  '
Set-TextValue "D13" '14.40'
$ws.Range("E13").Value = '  +0.25%  '
Set-TextValue "D14" '20.42'
$ws.Range("E14").Value = '  -0.95%  '
$ws.Range("E15").Value = '  -0.74%  '
$ws.Range("E16").Value = '  +0.99%  '
$ws.Range("D17").Value = '2.045.68'
$ws.Range("E17").Value = '  +0.99%  '
$ws.Range("D18").Value = '37.333.98'
$ws.Range("E18").Value = '  +1.05%  '
Set-TextValue "D19" '6.23'
$ws.Range("E19").Value = '  +1.21%  '
Set-TextValue "D20" '69.02'
$ws.Range("E20").Value = '  +0.15%  '
$ws.Range("D21").Value = '0.0₃0820'
$ws.Range("E21").Value = '  -0.74%  '
Set-TextValue "D22" '223.68'
$ws.Range("E22").Value = '  -1.42%  '
$ws.Range("E23").Value = '  -0.16%  '
$ws.Range("E24").Value = '  +1.91%  '
$ws.Range("E25").Value = '  -1.11%  '
Set-TextValue "D26" '164.78'
$ws.Range("E26").Value = '  -1.10%  '
$ws.Range("E27").Value = '  -1.12%  '
$ws.Range("E28").Value = '  +5.03%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("E30").Value = '  -1.45%  '
$ws.Range("E31").Value = '  +0.25%  '
Set-TextValue "D32" '4.48'
$ws.Range("E32").Value = '  -0.33%  '
Set-TextValue "D33" '0.0607'
$ws.Range("E33").Value = '  -1.48%  '
Set-TextValue "D34" '2.03'
$ws.Range("E34").Value = '  +10.91%  '
Set-TextValue "D35" '4.48'
$ws.Range("E35").Value = '  +0.52%  '
$ws.Range("E36").Value = '  -1.25%  '
Set-TextValue "D37" '3.23'
Set-TextValue "D38" '5.70'
$ws.Range("E38").Value = '  +6.37%  '
Set-TextValue "D39" '1.00'
$ws.Range("E39").Value = '  +0.25%  '
$ws.Range("D40").Value = '1.467.08'
$ws.Range("E40").Value = '  -1.62%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D41" '0.0214'
$ws.Range("E41").Value = '  -2.44%  '
$ws.Range("B42").Value = 'FTXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue "D42" '4.34'
$ws.Range("E42").Value = '  +20.33%  '
$ws.Range("E43").Value = '  +1.74%  '
Set-TextValue "D44" '0.0925'
$ws.Range("E44").Value = '  -0.60%  '
Set-TextValue "D45" '94.53'
$ws.Range("E45").Value = '  -0.94%  '
Set-TextValue "D46" '16.27'
$ws.Range("E46").Value = '  -5.12%  '
$ws.Range("E47").Value = '  -2.77%  '
$ws.Range("E48").Value = '  +0.36%  '
$ws.Range("E49").Value = '  -2.59%  '
Set-TextValue "D50" '2.94'
$ws.Range("E50").Value = '  +0.97%  '
$ws.Range("D51").Value = '2.221.16'
$ws.Range("E51").Value = '  +0.60%  '

Write-Host "Done applying cryptos update."
